$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 92, pushing the existing rows 92-152 down to 94-154.
$ws.Range("A92:A93").EntireRow.Insert()

# New row 92: Black Amber / Primera entry for 2022-01-13 (Región Metropolitana)
$ws.Cells.Item(92,1).Value = 8
$ws.Cells.Item(92,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(92,3).Value = "Coquimbo"
$ws.Cells.Item(92,4).Value = "2022-01-13"
$ws.Cells.Item(92,5).Value = 4
$ws.Cells.Item(92,6).Value = "Fruta"
$ws.Cells.Item(92,7).Value = 100103
$ws.Cells.Item(92,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(92,9).Value = 100103002
$ws.Cells.Item(92,10).Value = "Ciruela"
$ws.Cells.Item(92,11).Value = "Black Amber"
$ws.Cells.Item(92,12).Value = "Primera"
$ws.Cells.Item(92,13).Value = 20
$ws.Cells.Item(92,14).Value = 290000
$ws.Cells.Item(92,15).Value = 295000
$ws.Cells.Item(92,16).Value = 292500
$ws.Cells.Item(92,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(92,18).Value = "Región Metropolitana"
$ws.Cells.Item(92,19).Value = 650
$ws.Cells.Item(92,20).Value = 450

# New row 93: Black Amber / Segunda entry for 2022-01-13 (Región Metropolitana)
$ws.Cells.Item(93,1).Value = 8
$ws.Cells.Item(93,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(93,3).Value = "Coquimbo"
$ws.Cells.Item(93,4).Value = "2022-01-13"
$ws.Cells.Item(93,5).Value = 4
$ws.Cells.Item(93,6).Value = "Fruta"
$ws.Cells.Item(93,7).Value = 100103
$ws.Cells.Item(93,8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(93,9).Value = 100103002
$ws.Cells.Item(93,10).Value = "Ciruela"
$ws.Cells.Item(93,11).Value = "Black Amber"
$ws.Cells.Item(93,12).Value = "Segunda"
$ws.Cells.Item(93,13).Value = 20
$ws.Cells.Item(93,14).Value = 240000
$ws.Cells.Item(93,15).Value = 245000
$ws.Cells.Item(93,16).Value = 242500
$ws.Cells.Item(93,17).Value = "`$/bins (450 kilos)"
$ws.Cells.Item(93,18).Value = "Región Metropolitana"
$ws.Cells.Item(93,19).Value = 539
$ws.Cells.Item(93,20).Value = 450
